$d = $word.ActiveDocument

# --- 1. "Explain different types of storage structures..." paragraph ---
# The original is a single run; the target splits it into three runs
# (" Explain " / "d" / "ifferent types of storage structures ...") while
# keeping the visible text identical. Re-apply (identical) direct
# character formatting to just the single "d" so Word has to carve out a
# new run for it, without leaving any detectable formatting difference.
$rng = $d.Content
$rng.Find.Execute("different types of storage structures", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$charRng = $d.Range($rng.Start, $rng.Start + 1)
$charRng.Bold = 1
$charRng.Bold = 0
$charRng.Text = "d"

# --- 2. "Explain user mode and kernel mode of OS..how it works with diagram" ---
# Originally split across several runs (around proofing-error marks); the
# target collapses it back into a single run with the proofErr markers
# gone. A same-text Find/Replace normalizes the run & drops the stale
# proofErr bookkeeping.
$d.Content.Find.Execute("Explain user mode and kernel mode of OS..how it works with diagram", $true, $false, $false, $false, $false, $true, 1, $false, "Explain user mode and kernel mode of OS..how it works with diagram", 2)

# --- 3. "Explain different type of os structures like MS DOS, LINUX, Layerd, Microkernel with diagram" ---
$d.Content.Find.Execute("Explain different type of os structures like MS DOS, LINUX, Layerd, Microkernel with diagram", $true, $false, $false, $false, $false, $true, 1, $false, "Explain different type of os structures like MS DOS, LINUX, Layerd, Microkernel with diagram", 2)

# --- 4. "Explain difference among program, process and thread " ---
$d.Content.Find.Execute("Explain difference among program, process and thread ", $true, $false, $false, $false, $false, $true, 1, $false, "Explain difference among program, process and thread ", 2)

# --- 5. "Explain Process Scheduling with diagram(job, read, device queues)" ---
$d.Content.Find.Execute("Explain Process Scheduling with diagram(job, read, device queues)", $true, $false, $false, $false, $false, $true, 1, $false, "Explain Process Scheduling with diagram(job, read, device queues)", 2)

# --- 6a. "Blocking and Non blocking schemes" ---
$d.Content.Find.Execute("Blocking and Non blocking schemes", $true, $false, $false, $false, $false, $true, 1, $false, "Blocking and Non blocking schemes", 2)

# --- 6b. "Explain windows Lpc with diagram(last slide diagram) ***(most important)" ---
$d.Content.Find.Execute("Explain windows Lpc with diagram(last slide diagram) ***(most important)", $true, $false, $false, $false, $false, $true, 1, $false, "Explain windows Lpc with diagram(last slide diagram) ***(most important)", 2)
